# This workbook is a weekly price log. A new week's data is inserted at the
# top of the "Betarraga" data block (rows 300-301), and every subsequent
# week shifts down by one pair of rows (one row for "Primera" quality, one
# for "Segunda" quality). The two rows that fall off the bottom of the
# original range (old rows 410-411) become two brand-new rows (412-413)
# at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the existing data before any of it gets overwritten ---------
# Rows 300-409 (110 rows) will shift down to 302-411.
$blockMain = $ws.Range("A300:R409").Value2
# Rows 410-411 (last pair) will become the brand-new rows 412-413.
$blockTail = $ws.Range("A410:R411").Value2
# Remember the date-cell number format so the newly created rows 412-413
# keep the same formatting as the rest of column D.
$dateFormat = $ws.Range("D410").NumberFormat

# --- Extend the sheet with the two new rows at the bottom ----------------
$ws.Range("A412:R413").Value2 = $blockTail
$ws.Range("D412:D413").NumberFormat = $dateFormat

# --- Shift the existing 110 rows of data down by 2 rows -------------------
$ws.Range("A302:R411").Value2 = $blockMain

# --- Write the brand-new week of data into the freed-up rows 300-301 -----
# Row 300: "Primera" quality
$ws.Range("D300").Value2 = 44985
$ws.Range("J300").Value2 = 1800
$ws.Range("K300").Value2 = 500
$ws.Range("L300").Value2 = 600
$ws.Range("M300").Value2 = 550
$ws.Range("P300").Value2 = 183

# Row 301: "Segunda" quality
$ws.Range("D301").Value2 = 44985
$ws.Range("J301").Value2 = 1200
$ws.Range("K301").Value2 = 400
$ws.Range("L301").Value2 = 450
$ws.Range("M301").Value2 = 425
$ws.Range("P301").Value2 = 142
